$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1818255.9
$ws.Range("I5").Value = 2857225.5
$ws.Range("J5").Value = 58.75
$ws.Range("K5").Value = 2857225.5
$ws.Range("L5").Value = 58.75
$ws.Range("M5").Value = -2857110.5
$ws.Range("N5").Value = -288.75
$ws.Range("H28").Value = 14639.8
$ws.Range("I28").Value = 5787.5
$ws.Range("J28").Value = 24756.715
$ws.Range("K28").Value = 5787.5
$ws.Range("L28").Value = 24756.715
$ws.Range("M28").Value = -5302.5
$ws.Range("N28").Value = -25726.715
$ws.Range("H32").Value = 7771.4287
$ws.Range("J32").Value = 5128.5
$ws.Range("L32").Value = 5128.5
$ws.Range("N32").Value = -5780.5
$ws.Range("H51").Value = 7125.5
$ws.Range("J51").Value = 7125.5
$ws.Range("L51").Value = 7125.5
$ws.Range("N51").Value = -8093.5
$ws.Range("H53").Value = 1998.8
$ws.Range("J53").Value = 2373.5
$ws.Range("L53").Value = 2373.5
$ws.Range("N53").Value = -3647.5
$ws.Range("H113").Value = 5748.4736
$ws.Range("I113").Value = 7174.8887
$ws.Range("J113").Value = 4464.7
$ws.Range("K113").Value = 7174.8887
$ws.Range("L113").Value = 4464.7
$ws.Range("M113").Value = -3920.8887
$ws.Range("N113").Value = -10972.7
$ws.Range("H132").Value = 35150.4
$ws.Range("I132").Value = 1783.3334
$ws.Range("K132").Value = 5350.0002
$ws.Range("M132").Value = -2820.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 811.03705
$ws.Range("I2").Value = 773.1667
$ws.Range("K2").Value = 773.1667
$ws.Range("M2").Value = -660.1667
$ws.Range("H24").Value = 37749.5
$ws.Range("J24").Value = 37749.5
$ws.Range("L24").Value = 37749.5
$ws.Range("N24").Value = -38497.5
$ws.Range("H61").Value = 2839.12
$ws.Range("I61").Value = 2402.5
$ws.Range("K61").Value = 2402.5
$ws.Range("M61").Value = -2190.5
$ws.Range("H63").Value = 1853.44
$ws.Range("I63").Value = 1387.5264
$ws.Range("K63").Value = 1387.5264
$ws.Range("M63").Value = -701.5264
$ws.Range("H66").Value = 1853.44
$ws.Range("I66").Value = 1387.5264
$ws.Range("K66").Value = 6937.632
$ws.Range("M66").Value = -3505.632
$ws.Range("H100").Value = 37749.5
$ws.Range("J100").Value = 37749.5
$ws.Range("L100").Value = 37749.5
$ws.Range("N100").Value = -39913.5
$ws.Range("H110").Value = 2368.2
$ws.Range("I110").Value = 1775.8
$ws.Range("K110").Value = 1775.8
$ws.Range("M110").Value = 269.2
$ws.Range("H116").Value = 811.03705
$ws.Range("I116").Value = 773.1667
$ws.Range("K116").Value = 773.1667
$ws.Range("M116").Value = 1520.8333
$ws.Range("H122").Value = 1925.8182
$ws.Range("J122").Value = 2041.5714
$ws.Range("L122").Value = 6124.7142
$ws.Range("N122").Value = -11024.7142
$ws.Range("H136").Value = 2839.12
$ws.Range("I136").Value = 2402.5
$ws.Range("K136").Value = 7207.5
$ws.Range("M136").Value = -4657.5
$ws.Range("H138").Value = 33378330
$ws.Range("J138").Value = 67495
$ws.Range("L138").Value = 67495
$ws.Range("N138").Value = -77775

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 811.03705
$ws.Range("I3").Value = 773.1667
$ws.Range("K3").Value = 773.1667
$ws.Range("M3").Value = -659.1667
$ws.Range("H102").Value = 25587.143
$ws.Range("I102").Value = 25587.143
$ws.Range("K102").Value = 25587.143
$ws.Range("M102").Value = -22342.143
$ws.Range("H107").Value = 1562.1333
$ws.Range("I107").Value = 1143.909
$ws.Range("K107").Value = 1143.909
$ws.Range("M107").Value = 776.0909999999999
$ws.Range("H134").Value = 3895.1333
$ws.Range("I134").Value = 3823.2942
$ws.Range("K134").Value = 11469.8826
$ws.Range("M134").Value = -8934.882599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2082.7693
$ws.Range("I16").Value = 2038.8889
$ws.Range("K16").Value = 2038.8889
$ws.Range("M16").Value = -1751.8889
$ws.Range("H94").Value = 2394.3635
$ws.Range("I94").Value = 4080.8
$ws.Range("J94").Value = 989
$ws.Range("K94").Value = 4080.8
$ws.Range("L94").Value = 989
$ws.Range("M94").Value = -3629.8
$ws.Range("N94").Value = -1891
$ws.Range("H113").Value = 2082.7693
$ws.Range("I113").Value = 2038.8889
$ws.Range("K113").Value = 2038.8889
$ws.Range("M113").Value = 131.1111000000001
$ws.Range("H134").Value = 2872.3684
$ws.Range("I134").Value = 2798.8333
$ws.Range("K134").Value = 8396.499899999999
$ws.Range("M134").Value = -5861.499899999999
$ws.Range("H141").Value = 49994.6
$ws.Range("I141").Value = 24990
$ws.Range("K141").Value = 24990
$ws.Range("M141").Value = -19810

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 26
$ws.Range("I2").Value = 19.142857
$ws.Range("J2").Value = 32.857143
$ws.Range("K2").Value = 114.857142
$ws.Range("L2").Value = 197.142858
$ws.Range("M2").Value = -1.857141999999996
$ws.Range("N2").Value = -423.142858
$ws.Range("H4").Value = 10978055
$ws.Range("I4").Value = 2025212.2
$ws.Range("K4").Value = 6075636.6
$ws.Range("M4").Value = -6075524.6
$ws.Range("H6").Value = 77
$ws.Range("I6").Value = 77
$ws.Range("K6").Value = 231
$ws.Range("M6").Value = -118
$ws.Range("H7").Value = 338.125
$ws.Range("I7").Value = 388.33334
$ws.Range("K7").Value = 1165.00002
$ws.Range("M7").Value = -1053.00002
$ws.Range("H11").Value = 490
$ws.Range("I11").Value = 490
$ws.Range("K11").Value = 1470
$ws.Range("M11").Value = -1330
$ws.Range("H33").Value = 1223.625
$ws.Range("I33").Value = 118.75
$ws.Range("J33").Value = 2328.5
$ws.Range("K33").Value = 712.5
$ws.Range("L33").Value = 13971
$ws.Range("M33").Value = -429.5
$ws.Range("N33").Value = -14537
$ws.Range("H114").Value = 14856.571
$ws.Range("J114").Value = 16666.334
$ws.Range("L114").Value = 49999.00199999999
$ws.Range("N114").Value = -56507.00199999999
$ws.Range("H133").Value = 6181.727
$ws.Range("I133").Value = 999.5
$ws.Range("J133").Value = 7333.3335
$ws.Range("K133").Value = 2998.5
$ws.Range("L133").Value = 22000.0005
$ws.Range("M133").Value = 2061.5
$ws.Range("N133").Value = -32120.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 4003
$ws.Range("I24").Value = 5006
$ws.Range("J24").Value = 3000
$ws.Range("K24").Value = 5006
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = -4833
$ws.Range("N24").Value = -3346
$ws.Range("H70").Value = 9130.588
$ws.Range("I70").Value = 8286.1
$ws.Range("K70").Value = 8286.1
$ws.Range("M70").Value = -8016.1
$ws.Range("H73").Value = 9130.588
$ws.Range("I73").Value = 8286.1
$ws.Range("K73").Value = 8286.1
$ws.Range("M73").Value = -7350.1
$ws.Range("H132").Value = 3166.1924
$ws.Range("I132").Value = 3328.5264
$ws.Range("K132").Value = 9985.5792
$ws.Range("M132").Value = -7455.5792

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2151.3
$ws.Range("J46").Value = 2221.9412
$ws.Range("L46").Value = 2221.9412
$ws.Range("N46").Value = -2597.9412
$ws.Range("H55").Value = 217.14815
$ws.Range("I55").Value = 227.82353
$ws.Range("J55").Value = 199
$ws.Range("K55").Value = 227.82353
$ws.Range("L55").Value = 199
$ws.Range("M55").Value = -54.82353000000001
$ws.Range("N55").Value = -545
$ws.Range("H61").Value = 9992.571
$ws.Range("I61").Value = 9992
$ws.Range("J61").Value = 9994
$ws.Range("K61").Value = 9992
$ws.Range("L61").Value = 9994
$ws.Range("M61").Value = -9790
$ws.Range("N61").Value = -10398
$ws.Range("H113").Value = 9992.571
$ws.Range("I113").Value = 9992
$ws.Range("J113").Value = 9994
$ws.Range("K113").Value = 9992
$ws.Range("L113").Value = 9994
$ws.Range("M113").Value = -7822
$ws.Range("N113").Value = -14334
$ws.Range("H122").Value = 4298.4116
$ws.Range("I122").Value = 4208.4
$ws.Range("K122").Value = 12625.2
$ws.Range("M122").Value = -10175.2
$ws.Range("H132").Value = 3641.2917
$ws.Range("J132").Value = 3918.111
$ws.Range("L132").Value = 11754.333
$ws.Range("N132").Value = -16814.333
$ws.Range("H136").Value = 2768.5
$ws.Range("J136").Value = 2377.4
$ws.Range("L136").Value = 7132.200000000001
$ws.Range("N136").Value = -12232.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 20950
$ws.Range("I31").Value = 20800
$ws.Range("K31").Value = 20800
$ws.Range("M31").Value = -20452
$ws.Range("H113").Value = 913.4545000000001
$ws.Range("I113").Value = 579.6
$ws.Range("J113").Value = 1191.6666
$ws.Range("K113").Value = 1738.8
$ws.Range("L113").Value = 3574.9998
$ws.Range("M113").Value = 431.1999999999998
$ws.Range("N113").Value = -7914.9998
$ws.Range("H132").Value = 3687.0645
$ws.Range("I132").Value = 3644.9666
$ws.Range("K132").Value = 10934.8998
$ws.Range("M132").Value = -8404.899800000001
$ws.Range("H136").Value = 2608.8
$ws.Range("I136").Value = 2075.353
$ws.Range("K136").Value = 6226.059
$ws.Range("M136").Value = -3676.059
